$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case municipality/state name fixes (capitalize every word)
$ws.Range("B5").Value = "Pabellón De Arteaga"
$ws.Range("B6").Value = "Rincón De Romos"
$ws.Range("B18").Value = "Amatenango De La Frontera"
$ws.Range("B20").Value = "Bejucal De Ocampo"
$ws.Range("B22").Value = "Benemérito De Las Américas"
$ws.Range("B34").Value = "Comitán De Domínguez"
$ws.Range("B49").Value = "Mazapa De Madero"
$ws.Range("B58").Value = "San Cristóbal De Las Casas"
$ws.Range("B79").Value = "Hidalgo Del Parral"
$ws.Range("B97").Value = "Villa De Álvarez"
$ws.Range("A99").Value = "Ciudad De México"
$ws.Range("B103").Value = "Cuajimalpa De Morelos"
$ws.Range("B117").Value = "Coneto De Comonfort"
$ws.Range("B124").Value = "Pánuco De Coronado"
$ws.Range("B128").Value = "San Juan Del Río"
$ws.Range("A133").Value = "Estado De México"
$ws.Range("B133").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B136").Value = "Almoloya De Alquisiras"
$ws.Range("B147").Value = "Ecatepec De Morelos"
$ws.Range("B152").Value = "Ixtapan De La Sal"
$ws.Range("B162").Value = "Naucalpan De Juárez"
$ws.Range("B169").Value = "San Felipe Del Progreso"
$ws.Range("B177").Value = "Tenango Del Aire"
$ws.Range("B178").Value = "Tenango Del Valle"
$ws.Range("B184").Value = "Tlalnepantla De Baz"
$ws.Range("B189").Value = "Valle De Bravo"
$ws.Range("B190").Value = "Valle De Chalco Solidaridad"
$ws.Range("B191").Value = "Villa Del Carbón"
$ws.Range("B200").Value = "San Miguel De Allende"
$ws.Range("B201").Value = "Apaseo El Alto"
$ws.Range("B202").Value = "Apaseo El Grande"
$ws.Range("B209").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B212").Value = "Jaral Del Progreso"
$ws.Range("B220").Value = "Purísima Del Rincón"
$ws.Range("B225").Value = "San Francisco Del Rincón"
$ws.Range("B227").Value = "San Luis De La Paz"
$ws.Range("B228").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B230").Value = "Silao De La Victoria"
$ws.Range("B234").Value = "Valle De Santiago"
$ws.Range("B240").Value = "Acapulco De Juárez"
$ws.Range("B245").Value = "Atoyac De Álvarez"
$ws.Range("B246").Value = "Ayutla De Los Libres"
$ws.Range("B248").Value = "Buenavista De Cuéllar"
$ws.Range("B249").Value = "Chilapa De Álvarez"
$ws.Range("B250").Value = "Chilpancingo De Los Bravo"
$ws.Range("B252").Value = "Coyuca De Benítez"
$ws.Range("B253").Value = "Coyuca De Catalán"
$ws.Range("B257").Value = "Huitzuco De Los Figueroa"
$ws.Range("B258").Value = "Iguala De La Independencia"
$ws.Range("B259").Value = "Zihuatanejo De Azueta"
$ws.Range("B269").Value = "Taxco De Alarcón"
$ws.Range("B270").Value = "Técpan De Galeana"
$ws.Range("B272").Value = "Tixtla De Guerrero"
$ws.Range("B276").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B286").Value = "Atotonilco El Grande"
$ws.Range("B290").Value = "Cuautepec De Hinojosa"
$ws.Range("B296").Value = "Jacala De Ledezma"
$ws.Range("B300").Value = "Mineral Del Chico"
$ws.Range("B301").Value = "Mineral Del Monte"
$ws.Range("B302").Value = "Mixquiahuala De Juárez"
$ws.Range("B304").Value = "Pachuca De Soto"
$ws.Range("B309").Value = "Tenango De Doria"
$ws.Range("B310").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B311").Value = "Tezontepec De Aldama"
$ws.Range("B315").Value = "Tula De Allende"
$ws.Range("B316").Value = "Tulancingo De Bravo"
$ws.Range("B318").Value = "Zacualtipán De Ángeles"
$ws.Range("B323").Value = "Atemajac De Brizuela"
$ws.Range("B324").Value = "Autlán De Navarro"
$ws.Range("B330").Value = "Encarnación De Díaz"
$ws.Range("B334").Value = "Jilotlán De Los Dolores"
$ws.Range("B336").Value = "Lagos De Moreno"
$ws.Range("B340").Value = "Tamazula De Gordiano"
$ws.Range("B343").Value = "Tizapán El Alto"
$ws.Range("B344").Value = "Tlajomulco De Zúñiga"
$ws.Range("B349").Value = "Unión De San Antonio"
$ws.Range("B350").Value = "Unión De Tula"
$ws.Range("B400").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B423").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B426").Value = "Puente De Ixtla"
$ws.Range("B428").Value = "Tlaltizapán De Zapata"
$ws.Range("B436").Value = "Amatlán De Cañas"
$ws.Range("B442").Value = "Santa María Del Oro"
$ws.Range("B454").Value = "San Nicolás De Los Garza"
$ws.Range("B457").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B461").Value = "Ciénega De Zimatlán"
$ws.Range("B463").Value = "Fresnillo De Trujano"
$ws.Range("B464").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B465").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B466").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B467").Value = "Ixtlán De Juárez"
$ws.Range("B468").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B476").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B479").Value = "Oaxaca De Juárez"
$ws.Range("B480").Value = "Ocotlán De Morelos"
$ws.Range("B481").Value = "Pinotepa De Don Luis"
$ws.Range("B483").Value = "Putla Villa De Guerrero"
$ws.Range("B488").Value = "San Antonino El Alto"
$ws.Range("B490").Value = "San Antonio De La Cal"
$ws.Range("B521").Value = "San Mateo Del Mar"
$ws.Range("B525").Value = "San Miguel Del Puerto"
$ws.Range("B551").Value = "Santa Inés Del Monte"
$ws.Range("B581").Value = "Santo Domingo De Morelos"
$ws.Range("B587").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B588").Value = "Tanetze De Zaragoza"
$ws.Range("B589").Value = "Teotitlán De Flores Magón"
$ws.Range("B590").Value = "Tepelmeme Villa De Morelos"
$ws.Range("B591").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B592").Value = "Tlacolula De Matamoros"
$ws.Range("B593").Value = "Tlalixtac De Cabrera"
$ws.Range("B594").Value = "Totontepec Villa De Morelos"
$ws.Range("B595").Value = "Villa De Chilapa De Díaz"
$ws.Range("B596").Value = "Villa De Etla"
$ws.Range("B597").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B598").Value = "Villa De Zaachila"
$ws.Range("B600").Value = "Villa Sola De Vega"
$ws.Range("B601").Value = "Zimatlán De Álvarez"
$ws.Range("B607").Value = "Ayotoxco De Guerrero"
$ws.Range("B608").Value = "Chalchicomula De Sesma"
$ws.Range("B624").Value = "Ixcamilpa De Guerrero"
$ws.Range("B626").Value = "Izúcar De Matamoros"
$ws.Range("B643").Value = "Tecali De Herrera"
$ws.Range("B647").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B651").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B657").Value = "Xayacatlán De Bravo"
$ws.Range("B667").Value = "Amealco De Bonfil"
$ws.Range("B672").Value = "Jalpan De Serra"
$ws.Range("B673").Value = "Landa De Matamoros"
$ws.Range("B674").Value = "Pinal De Amoles"
$ws.Range("B676").Value = "San Juan Del Río"
$ws.Range("B689").Value = "Ciudad Del Maíz"
$ws.Range("B697").Value = "San Ciro De Acosta"
$ws.Range("B699").Value = "Santa María Del Río"
$ws.Range("B703").Value = "Tanquián De Escobedo"
$ws.Range("B705").Value = "Villa De Guadalupe"
$ws.Range("B706").Value = "Villa De Ramos"
$ws.Range("B707").Value = "Villa De Reyes"
$ws.Range("B753").Value = "Soto La Marina"
$ws.Range("B762").Value = "Contla De Juan Cuamatzi"
$ws.Range("B764").Value = "San Pablo Del Monte"
$ws.Range("B765").Value = "Sanctórum De Lázaro Cárdenas"
$ws.Range("B767").Value = "Tetla De La Solidaridad"
$ws.Range("B770").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B772").Value = "Amatlán De Los Reyes"
$ws.Range("B780").Value = "Cazones De Herrera"
$ws.Range("B790").Value = "Cosamaloapan De Carpio"
$ws.Range("B791").Value = "Cosautlán De Carvajal"
$ws.Range("B801").Value = "Ignacio De La Llave"
$ws.Range("B802").Value = "Ixhuatlán De Madero"
$ws.Range("B803").Value = "Ixhuatlán Del Café"
$ws.Range("B806").Value = "Juchique De Ferrer"
$ws.Range("B809").Value = "Martínez De La Torre"
$ws.Range("B818").Value = "Paso De Ovejas"
$ws.Range("B820").Value = "Poza Rica De Hidalgo"
$ws.Range("B826").Value = "Sayula De Alemán"
$ws.Range("B827").Value = "Soledad De Doblado"
$ws.Range("B837").Value = "Tlacotepec De Mejía"
$ws.Range("B843").Value = "Vega De Alatorre"
$ws.Range("B849").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B852").Value = "Cañitas De Felipe Pescador"
$ws.Range("B855").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B860").Value = "Jiménez Del Teul"
$ws.Range("B869").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B870").Value = "Trinidad García De La Cadena"
$ws.Range("B872").Value = "Villa De Cos"

# 3. Fix floating point rounding for percentage values where n=4
$ws.Range("D3").Value = 0.0009720534629404616
$ws.Range("D8").Value = 0.0009720534629404616
$ws.Range("D13").Value = 0.0009720534629404616
$ws.Range("D14").Value = 0.0009720534629404616
$ws.Range("D122").Value = 0.0009720534629404616
$ws.Range("D128").Value = 0.0009720534629404616
$ws.Range("D131").Value = 0.0009720534629404616
$ws.Range("D140").Value = 0.0009720534629404616
$ws.Range("D141").Value = 0.0009720534629404616
$ws.Range("D152").Value = 0.0009720534629404616
$ws.Range("D158").Value = 0.0009720534629404616
$ws.Range("D171").Value = 0.0009720534629404616
$ws.Range("D174").Value = 0.0009720534629404616
$ws.Range("D200").Value = 0.0009720534629404616
$ws.Range("D208").Value = 0.0009720534629404616
$ws.Range("D220").Value = 0.0009720534629404616
$ws.Range("D227").Value = 0.0009720534629404616
$ws.Range("D244").Value = 0.0009720534629404616
$ws.Range("D270").Value = 0.0009720534629404616
$ws.Range("D273").Value = 0.0009720534629404616
$ws.Range("D300").Value = 0.0009720534629404616
$ws.Range("D356").Value = 0.0009720534629404616
$ws.Range("D360").Value = 0.0009720534629404616
$ws.Range("D363").Value = 0.0009720534629404616
$ws.Range("D384").Value = 0.0009720534629404616
$ws.Range("D412").Value = 0.0009720534629404616
$ws.Range("D416").Value = 0.0009720534629404616
$ws.Range("D420").Value = 0.0009720534629404616
$ws.Range("D425").Value = 0.0009720534629404616
$ws.Range("D433").Value = 0.0009720534629404616
$ws.Range("D474").Value = 0.0009720534629404616
$ws.Range("D489").Value = 0.0009720534629404616
$ws.Range("D491").Value = 0.0009720534629404616
$ws.Range("D519").Value = 0.0009720534629404616
$ws.Range("D524").Value = 0.0009720534629404616
$ws.Range("D557").Value = 0.0009720534629404616
$ws.Range("D575").Value = 0.0009720534629404616
$ws.Range("D582").Value = 0.0009720534629404616
$ws.Range("D592").Value = 0.0009720534629404616
$ws.Range("D667").Value = 0.0009720534629404616
$ws.Range("D710").Value = 0.0009720534629404616
$ws.Range("D723").Value = 0.0009720534629404616
$ws.Range("D734").Value = 0.0009720534629404616
$ws.Range("D760").Value = 0.0009720534629404616
$ws.Range("D775").Value = 0.0009720534629404616
$ws.Range("D782").Value = 0.0009720534629404616
$ws.Range("D793").Value = 0.0009720534629404616
$ws.Range("D814").Value = 0.0009720534629404616
$ws.Range("D819").Value = 0.0009720534629404616
$ws.Range("D820").Value = 0.0009720534629404616
$ws.Range("D824").Value = 0.0009720534629404616
$ws.Range("D843").Value = 0.0009720534629404616
$ws.Range("D857").Value = 0.0009720534629404616

# 4. Remove trailing metadata/footer rows (879-883) and shrink dimension accordingly
$ws.Rows("879:883").Delete()
